# On the "Macros" slide, the paragraph "so ext-if=<smart-quote>em0<smart-quote> will not work"
# is currently split into two runs: '-if="em0" will ' and 'not work'.
# Merge them into a single run reading '-if="em0" will not work' (keeping the
# formatting of the first of the two runs), leaving the surrounding runs
# ("so " and "ext") untouched.
#
# Note: TextRange.Text reports the curly ("smart") quotes already stored in
# the deck as plain ASCII quotes, so the search needle below uses a plain
# quote while the replacement text writes back the curly quote characters
# that belong in the OOXML.

$p = $ppt.ActivePresentation

$quoteOpen  = [char]0x201C
$quoteClose = [char]0x201D
$needle     = '-if="em0" will '
$newText    = '-if=' + $quoteOpen + 'em0' + $quoteClose + ' will not work'

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if (-not $shape.HasTextFrame) { continue }
        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf($needle)
        if ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, $newText.Length)
            $sub.Text = $newText
        }
    }
}
